$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44497
$ws.Range("J2").Value = 20
$ws.Range("K2").Value = 4000
$ws.Range("L2").Value = 4000
$ws.Range("M2").Value = 4000
$ws.Range("P2").Value = 4000

# Row 3
$ws.Range("D3").Value = 44259
$ws.Range("J3").Value = 30

# Row 4
$ws.Range("D4").Value = 44301
$ws.Range("J4").Value = 40
$ws.Range("K4").Value = 3000
$ws.Range("L4").Value = 3000
$ws.Range("M4").Value = 3000
$ws.Range("P4").Value = 3000

# Row 5
$ws.Range("D5").Value = 44679
$ws.Range("J5").Value = 50
$ws.Range("K5").Value = 5000
$ws.Range("L5").Value = 5000
$ws.Range("M5").Value = 5000
$ws.Range("P5").Value = 5000

# Row 6
$ws.Range("D6").Value = 44680
$ws.Range("J6").Value = 20
$ws.Range("K6").Value = 5000
$ws.Range("L6").Value = 5000
$ws.Range("M6").Value = 5000
$ws.Range("P6").Value = 5000

# Row 7 unchanged

# Row 8
$ws.Range("D8").Value = 44498
$ws.Range("J8").Value = 40

# Row 9
$ws.Range("D9").Value = 44176
$ws.Range("J9").Value = 10

# Row 10
$ws.Range("D10").Value = 44656
$ws.Range("J10").Value = 85

# Row 11
$ws.Range("D11").Value = 44312
$ws.Range("J11").Value = 50
$ws.Range("K11").Value = 4000
$ws.Range("L11").Value = 4000
$ws.Range("M11").Value = 4000
$ws.Range("P11").Value = 4000

# Row 12
$ws.Range("D12").Value = 44390
$ws.Range("J12").Value = 55
$ws.Range("K12").Value = 6000
$ws.Range("L12").Value = 6000
$ws.Range("M12").Value = 6000
$ws.Range("P12").Value = 6000

# Row 13
$ws.Range("D13").Value = 44508
$ws.Range("J13").Value = 30
$ws.Range("K13").Value = 4000
$ws.Range("L13").Value = 4000
$ws.Range("M13").Value = 4000
$ws.Range("P13").Value = 4000

# Row 14
$ws.Range("D14").Value = 44313
$ws.Range("J14").Value = 20

# Row 15
$ws.Range("D15").Value = 44291
$ws.Range("J15").Value = 35

# Row 16
$ws.Range("D16").Value = 44649
$ws.Range("J16").Value = 20
$ws.Range("K16").Value = 5000
$ws.Range("L16").Value = 5000
$ws.Range("M16").Value = 5000
$ws.Range("P16").Value = 5000

# Row 17
$ws.Range("D17").Value = 44280
$ws.Range("J17").Value = 55

# Row 18
$ws.Range("D18").Value = 44315
$ws.Range("K18").Value = 4000
$ws.Range("L18").Value = 4000
$ws.Range("M18").Value = 4000
$ws.Range("P18").Value = 4000

# Row 19
$ws.Range("D19").Value = 44316
$ws.Range("J19").Value = 20

# Row 20
$ws.Range("D20").Value = 44504
$ws.Range("J20").Value = 55

# Row 21
$ws.Range("D21").Value = 44365
$ws.Range("J21").Value = 55
$ws.Range("K21").Value = 5000
$ws.Range("L21").Value = 5000
$ws.Range("M21").Value = 5000
$ws.Range("P21").Value = 5000
